# Update requirements and action items
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Lit. Rev. item text for several rows to reflect the new "SLAM" literature review topics
$ws.Range("F10").Value = "Lit. Rev: SLAM"
$ws.Range("F9").Value = "Lit. Rev.: SLAM"
$ws.Range("F11").Value = "Lit. Rev.: SLAM"
$ws.Range("F6").Value = "Lit. Rev.: Visual Odemetry / SLAM"

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("F7").Select()
